$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.810.13'
$ws.Range("E2").Value = '  +1.16%  '
$ws.Range("D3").Value = '3.458.94'
$ws.Range("E3").Value = '  +2.80%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '580.34'
$ws.Range("E5").Value = '  +1.71%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '148.52'
$ws.Range("E6").Value = '  +9.54%  '
$ws.Range("D7").Value = '3.459.76'
$ws.Range("E7").Value = '  +2.87%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("E9").Value = '  +1.20%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.76'
$ws.Range("E10").Value = '  +3.57%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.125'
$ws.Range("E11").Value = '  +1.55%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.391'
$ws.Range("E12").Value = '  +1.84%  '
$ws.Range("D13").Value = '4.050.57'
$ws.Range("E13").Value = '  +2.89%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.23'
$ws.Range("E14").Value = '  +8.84%  '
$ws.Range("E15").Value = '  -0.34%  '
$ws.Range("E16").Value = '  +1.73%  '
$ws.Range("D17").Value = '3.444.67'
$ws.Range("E17").Value = '  +2.63%  '
$ws.Range("D18").Value = '61.859.08'
$ws.Range("E18").Value = '  +1.16%  '
$ws.Range("E19").Value = '  +9.36%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.32'
$ws.Range("E20").Value = '  +2.40%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.45'
$ws.Range("E21").Value = '  +2.27%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '385.52'
$ws.Range("E22").Value = '  +2.38%  '
$ws.Range("E23").Value = '  +3.26%  '
$ws.Range("D24").Value = '3.597.47'
$ws.Range("E24").Value = '  +2.98%  '
$ws.Range("E25").Value = '  +0.20%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '72.61'
$ws.Range("E27").Value = '  +2.23%  '
$ws.Range("E28").Value = '  -1.86%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.180'
$ws.Range("E29").Value = '  +8.93%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.82'
$ws.Range("E30").Value = '  +5.15%  '
$ws.Range("B31").Value = 'Binance-PegBSC-USD'
$ws.Range("C31").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.999'
$ws.Range("E31").Value = '  -0.50%  '
$ws.Range("B32").Value = 'Fetch.AI'
$ws.Range("C32").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.53'
$ws.Range("E32").Value = '  -13.11%  '
$ws.Range("E33").Value = '  +1.41%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.17'
$ws.Range("E34").Value = '  +1.87%  '
$ws.Range("E35").Value = '  +0.00%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '23.96'
$ws.Range("E36").Value = '  +1.93%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '7.07'
$ws.Range("E37").Value = '  +4.65%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.21'
$ws.Range("E38").Value = '  +0.44%  '
$ws.Range("E39").Value = '  +2.75%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '166.52'
$ws.Range("E40").Value = '  +0.96%  '
$ws.Range("E41").Value = '  +4.78%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.797'
$ws.Range("E42").Value = '  +3.70%  '
$ws.Range("E43").Value = '  +9.73%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.73'
$ws.Range("E44").Value = '  +1.72%  '
$ws.Range("E45").Value = '  +0.09%  '
$ws.Range("B46").Value = 'OKB'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '42.36'
$ws.Range("E46").Value = '  +2.17%  '
$ws.Range("B47").Value = 'Filecoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.48'
$ws.Range("E47").Value = '  +2.23%  '
$ws.Range("E48").Value = '  -2.26%  '
$ws.Range("D49").Value = '2.608.85'
$ws.Range("E49").Value = '  +11.43%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.96'
$ws.Range("E50").Value = '  +2.42%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '23.35'
$ws.Range("E51").Value = '  +0.52%  '
